$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted numbers (e.g. "52.081.57", "351.81").
# Force the whole data range to Text format first so Excel does not
# auto-convert the assigned strings into numeric values, then restore
# the default "Normal" style so the cell styling is unchanged.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '52.081.57'
$ws.Range("D3").Value = '2.894.15'
$ws.Range("E3").Value = '  +3.39%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '351.81'
$ws.Range("E5").Value = '  -0.80%  '
$ws.Range("D6").Value = '111.52'
$ws.Range("E6").Value = '  +2.11%  '
$ws.Range("D7").Value = '0.559'
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '39.87'
$ws.Range("E10").Value = '  -0.37%  '
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").Value = '0.0857'
$ws.Range("E12").Value = '  +2.28%  '
$ws.Range("D13").Value = '19.96'
$ws.Range("E13").Value = '  -0.14%  '
$ws.Range("D14").Value = '7.78'
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").Value = '3.347.76'
$ws.Range("E15").Value = '  +3.31%  '
$ws.Range("E16").Value = '  +7.20%  '
$ws.Range("D17").Value = '2.906.11'
$ws.Range("E17").Value = '  +3.78%  '
$ws.Range("D18").Value = '52.072.18'
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("D19").Value = '7.72'
$ws.Range("E19").Value = '  -0.26%  '
$ws.Range("E20").Value = '  +5.94%  '
$ws.Range("D21").Value = '14.33'
$ws.Range("E21").Value = '  +7.60%  '
$ws.Range("D22").Value = '0.0₃0980'
$ws.Range("E22").Value = '  +0.80%  '
$ws.Range("E23").Value = '  +0.29%  '
$ws.Range("D24").Value = '269.43'
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("D26").Value = '26.51'
$ws.Range("E26").Value = '  +1.89%  '
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("D28").Value = '0.164'
$ws.Range("E28").Value = '  -0.17%  '
$ws.Range("D29").Value = '10.52'
$ws.Range("E29").Value = '  +1.67%  '
$ws.Range("D30").Value = '38.42'
$ws.Range("E30").Value = '  +2.61%  '
$ws.Range("E31").Value = '  +0.71%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '6.43'
$ws.Range("E32").Value = '  +3.40%  '
$ws.Range("B33").Value = 'RenderToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D33").Value = '6.15'
$ws.Range("E33").Value = '  +8.76%  '
$ws.Range("D34").Value = '0.0947'
$ws.Range("E34").Value = '  +10.84%  '
$ws.Range("D35").Value = '53.02'
$ws.Range("E35").Value = '  +1.48%  '
$ws.Range("D36").Value = '0.0458'
$ws.Range("E36").Value = '  +2.55%  '
$ws.Range("E37").Value = '  -0.20%  '
$ws.Range("D38").Value = '3.31'
$ws.Range("E38").Value = '  +5.22%  '
$ws.Range("D39").Value = '18.64'
$ws.Range("E39").Value = '  -0.44%  '
$ws.Range("E40").Value = '  +2.79%  '
$ws.Range("E41").Value = '  +6.40%  '
$ws.Range("E42").Value = '  +1.66%  '
$ws.Range("D43").Value = '22.72'
$ws.Range("E43").Value = '  +3.81%  '
$ws.Range("D44").Value = '122.12'
$ws.Range("E44").Value = '  +1.94%  '
$ws.Range("E45").Value = '  +0.61%  '
$ws.Range("D46").Value = '3.58'
$ws.Range("E46").Value = '  +4.23%  '
$ws.Range("D47").Value = '2.200.93'
$ws.Range("E47").Value = '  +2.85%  '
$ws.Range("E48").Value = '  +6.21%  '
$ws.Range("E49").Value = '  +22.13%  '
$ws.Range("D50").Value = '0.948'
$ws.Range("E50").Value = '  +2.56%  '
$ws.Range("D51").Value = '0.0324'
$ws.Range("E51").Value = '  +11.37%  '

$ws.Range("D2:D51").Style = "Normal"
